# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E3) and
# "Correspond Handback DateTime" (H3) timestamps for the
# 0d24fb0b-0e88-4076-aee9-74cd73159bad handback row on both the
# zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E3").Value = "2016-03-14 03:31:58"
$zhcn.Range("H3").Value = "2016-03-14 03:32:17"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E3").Value = "2016-03-14 03:32:01"
$dede.Range("H3").Value = "2016-03-14 03:32:22"
